# "Linking leadership team to events"
# Repoint the three Speaker(s) cells that referenced Lauren Chenarides,
# Drew Hanks and Andi Carlson's old personal/department bio pages so they
# instead link to the project's own Leadership_team page.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Day 1 Opening Remarks) and Row 10 (Day 1 Wrap-up) both list
# Lauren Chenarides & Drew Hanks as speakers.
$laurenDrew = "[Lauren Chenarides](https://dataifa.github.io/difa-project/Leadership_team.html), [Drew Hanks](https://dataifa.github.io/difa-project/Leadership_team.html)"
$ws.Range("F2").Value = $laurenDrew
$ws.Range("F10").Value = $laurenDrew

# Row 8 (Session 3 / Developments in data linkages) lists Andi Carlson.
$ws.Range("F8").Value = "[Andi Carlson](https://dataifa.github.io/difa-project/Leadership_team.html)"

# Match the author's last active selection in the sheet.
$ws.Range("E16").Select()
